# Update NATMI LR-pair statistics for Thbs1-Itgb3 (rows 2-17, columns E-T)
# following re-computation with additional replicate ("Dr Hou advice").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (E..T): E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$data = New-Object 'object[,]' 16,16
# Row 2
$data[0,0] = 3.0  # E2
$data[0,1] = 1.0  # F2
$data[0,2] = 18.82499266666667  # G2
$data[0,3] = 56.474978  # H2
$data[0,4] = 0.06886869772378311  # I2
$data[0,5] = 0.0688686977237831  # J2
$data[0,6] = 3.0  # K2
$data[0,7] = 1.0  # L2
$data[0,8] = 3.778439  # M2
$data[0,9] = 11.335317  # N2
$data[0,10] = 0.4252971528324392  # O2
$data[0,11] = 0.4252971528324392  # P2
$data[0,12] = 71.12908646644733  # Q2
$data[0,13] = 640.161778198026  # R2
$data[0,14] = 0.02928966106120285  # S2
$data[0,15] = 0.02928966106120284  # T2
# Row 3
$data[1,0] = 3.0  # E3
$data[1,1] = 1.0  # F3
$data[1,2] = 18.82499266666667  # G3
$data[1,3] = 56.474978  # H3
$data[1,4] = 0.06886869772378311  # I3
$data[1,5] = 0.0688686977237831  # J3
$data[1,6] = 3.0  # K3
$data[1,7] = 1.0  # L3
$data[1,8] = 4.333403333333333  # M3
$data[1,9] = 13.00021  # N3
$data[1,10] = 0.4877633593505858  # O3
$data[1,11] = 0.4877633593505858  # P3
$data[1,12] = 81.57628597170888  # Q3
$data[1,13] = 734.18657374538  # R3
$data[1,14] = 0.03359162735585249  # S3
$data[1,15] = 0.03359162735585249  # T3
# Row 4
$data[2,0] = 3.0  # E4
$data[2,1] = 1.0  # F4
$data[2,2] = 18.82499266666667  # G4
$data[2,3] = 56.474978  # H4
$data[2,4] = 0.06886869772378311  # I4
$data[2,5] = 0.0688686977237831  # J4
$data[2,6] = 3.0  # K4
$data[2,7] = 1.0  # L4
$data[2,8] = 0.2909853333333334  # M4
$data[2,9] = 0.8729560000000001  # N4
$data[2,10] = 0.03275300561492853  # O4
$data[2,11] = 0.03275300561492853  # P4
$data[2,12] = 5.477796766107557  # Q4
$data[2,13] = 49.30017089496801  # R4
$data[2,14] = 0.002255656843239884  # S4
$data[2,15] = 0.002255656843239884  # T4
# Row 5
$data[3,0] = 3.0  # E5
$data[3,1] = 1.0  # F5
$data[3,2] = 18.82499266666667  # G5
$data[3,3] = 56.474978  # H5
$data[3,4] = 0.06886869772378311  # I5
$data[3,5] = 0.0688686977237831  # J5
$data[3,6] = 3.0  # K5
$data[3,7] = 1.0  # L5
$data[3,8] = 0.4814053333333333  # M5
$data[3,9] = 1.444216  # N5
$data[3,10] = 0.0541864822020464  # O5
$data[3,11] = 0.05418648220204641  # P5
$data[3,12] = 9.062451869694222  # Q5
$data[3,13] = 81.562066827248  # R5
$data[3,14] = 0.003731752463487887  # S5
$data[3,15] = 0.003731752463487887  # T5
# Row 6
$data[4,0] = 3.0  # E6
$data[4,1] = 1.0  # F6
$data[4,2] = 121.8208923333333  # G6
$data[4,3] = 365.462677  # H6
$data[4,4] = 0.4456653109566078  # I6
$data[4,5] = 0.4456653109566078  # J6
$data[4,6] = 3.0  # K6
$data[4,7] = 1.0  # L6
$data[4,8] = 3.778439  # M6
$data[4,9] = 11.335317  # N6
$data[4,10] = 0.4252971528324392  # O6
$data[4,11] = 0.4252971528324392  # P6
$data[4,12] = 460.2928106070677  # Q6
$data[4,13] = 4142.63529546361  # R6
$data[4,14] = 0.189540187866029  # S6
$data[4,15] = 0.189540187866029  # T6
# Row 7
$data[5,0] = 3.0  # E7
$data[5,1] = 1.0  # F7
$data[5,2] = 121.8208923333333  # G7
$data[5,3] = 365.462677  # H7
$data[5,4] = 0.4456653109566078  # I7
$data[5,5] = 0.4456653109566078  # J7
$data[5,6] = 3.0  # K7
$data[5,7] = 1.0  # L7
$data[5,8] = 4.333403333333333  # M7
$data[5,9] = 13.00021  # N7
$data[5,10] = 0.4877633593505858  # O7
$data[5,11] = 0.4877633593505858  # P7
$data[5,12] = 527.8990609069078  # Q7
$data[5,13] = 4751.091548162171  # R7
$data[5,14] = 0.2173792092182185  # S7
$data[5,15] = 0.2173792092182185  # T7
# Row 8
$data[6,0] = 3.0  # E8
$data[6,1] = 1.0  # F8
$data[6,2] = 121.8208923333333  # G8
$data[6,3] = 365.462677  # H8
$data[6,4] = 0.4456653109566078  # I8
$data[6,5] = 0.4456653109566078  # J8
$data[6,6] = 3.0  # K8
$data[6,7] = 1.0  # L8
$data[6,8] = 0.2909853333333334  # M8
$data[6,9] = 0.8729560000000001  # N8
$data[6,10] = 0.03275300561492853  # O8
$data[6,11] = 0.03275300561492853  # P8
$data[6,12] = 35.44809296257912  # Q8
$data[6,13] = 319.032836663212  # R8
$data[6,14] = 0.01459687843214065  # S8
$data[6,15] = 0.01459687843214065  # T8
# Row 9
$data[7,0] = 3.0  # E9
$data[7,1] = 1.0  # F9
$data[7,2] = 121.8208923333333  # G9
$data[7,3] = 365.462677  # H9
$data[7,4] = 0.4456653109566078  # I9
$data[7,5] = 0.4456653109566078  # J9
$data[7,6] = 3.0  # K9
$data[7,7] = 1.0  # L9
$data[7,8] = 0.4814053333333333  # M9
$data[7,9] = 1.444216  # N9
$data[7,10] = 0.0541864822020464  # O9
$data[7,11] = 0.05418648220204641  # P9
$data[7,12] = 58.64522728069245  # Q9
$data[7,13] = 527.8070455262321  # R9
$data[7,14] = 0.0241490354402197  # S9
$data[7,15] = 0.0241490354402197  # T9
# Row 10
$data[8,0] = 3.0  # E10
$data[8,1] = 1.0  # F10
$data[8,2] = 87.673585  # G10
$data[8,3] = 263.020755  # H10
$data[8,4] = 0.3207419907481189  # I10
$data[8,5] = 0.3207419907481188  # J10
$data[8,6] = 3.0  # K10
$data[8,7] = 1.0  # L10
$data[8,8] = 3.778439  # M10
$data[8,9] = 11.335317  # N10
$data[8,10] = 0.4252971528324392  # O10
$data[8,11] = 0.4252971528324392  # P10
$data[8,12] = 331.269292833815  # Q10
$data[8,13] = 2981.423635504335  # R10
$data[8,14] = 0.1364106554589835  # S10
$data[8,15] = 0.1364106554589835  # T10
# Row 11
$data[9,0] = 3.0  # E11
$data[9,1] = 1.0  # F11
$data[9,2] = 87.673585  # G11
$data[9,3] = 263.020755  # H11
$data[9,4] = 0.3207419907481189  # I11
$data[9,5] = 0.3207419907481188  # J11
$data[9,6] = 3.0  # K11
$data[9,7] = 1.0  # L11
$data[9,8] = 4.333403333333333  # M11
$data[9,9] = 13.00021  # N11
$data[9,10] = 0.4877633593505858  # O11
$data[9,11] = 0.4877633593505858  # P11
$data[9,12] = 379.9250054842833  # Q11
$data[9,13] = 3419.32504935855  # R11
$data[9,14] = 0.156446190892097  # S11
$data[9,15] = 0.156446190892097  # T11
# Row 12
$data[10,0] = 3.0  # E12
$data[10,1] = 1.0  # F12
$data[10,2] = 87.673585  # G12
$data[10,3] = 263.020755  # H12
$data[10,4] = 0.3207419907481189  # I12
$data[10,5] = 0.3207419907481188  # J12
$data[10,6] = 3.0  # K12
$data[10,7] = 1.0  # L12
$data[10,8] = 0.2909853333333334  # M12
$data[10,9] = 0.8729560000000001  # N12
$data[10,10] = 0.03275300561492853  # O12
$data[10,11] = 0.03275300561492853  # P12
$data[10,12] = 25.51172735575334  # Q12
$data[10,13] = 229.60554620178  # R12
$data[10,14] = 0.01050526422391649  # S12
$data[10,15] = 0.01050526422391649  # T12
# Row 13
$data[11,0] = 3.0  # E13
$data[11,1] = 1.0  # F13
$data[11,2] = 87.673585  # G13
$data[11,3] = 263.020755  # H13
$data[11,4] = 0.3207419907481189  # I13
$data[11,5] = 0.3207419907481188  # J13
$data[11,6] = 3.0  # K13
$data[11,7] = 1.0  # L13
$data[11,8] = 0.4814053333333333  # M13
$data[11,9] = 1.444216  # N13
$data[11,10] = 0.0541864822020464  # O13
$data[11,11] = 0.05418648220204641  # P13
$data[11,12] = 42.20653141145333  # Q13
$data[11,13] = 379.85878270308  # R13
$data[11,14] = 0.01737988017312187  # S13
$data[11,15] = 0.01737988017312187  # T13
# Row 14
$data[12,0] = 3.0  # E14
$data[12,1] = 1.0  # F14
$data[12,2] = 45.02666966666666  # G14
$data[12,3] = 135.080009  # H14
$data[12,4] = 0.1647240005714903  # I14
$data[12,5] = 0.1647240005714903  # J14
$data[12,6] = 3.0  # K14
$data[12,7] = 1.0  # L14
$data[12,8] = 3.778439  # M14
$data[12,9] = 11.335317  # N14
$data[12,10] = 0.4252971528324392  # O14
$data[12,11] = 0.4252971528324392  # P14
$data[12,12] = 170.1305247086503  # Q14
$data[12,13] = 1531.174722377853  # R14
$data[12,14] = 0.0700566484462239  # S14
$data[12,15] = 0.0700566484462239  # T14
# Row 15
$data[13,0] = 3.0  # E15
$data[13,1] = 1.0  # F15
$data[13,2] = 45.02666966666666  # G15
$data[13,3] = 135.080009  # H15
$data[13,4] = 0.1647240005714903  # I15
$data[13,5] = 0.1647240005714903  # J15
$data[13,6] = 3.0  # K15
$data[13,7] = 1.0  # L15
$data[13,8] = 4.333403333333333  # M15
$data[13,9] = 13.00021  # N15
$data[13,10] = 0.4877633593505858  # O15
$data[13,11] = 0.4877633593505858  # P15
$data[13,12] = 195.1187204224322  # Q15
$data[13,13] = 1756.06848380189  # R15
$data[13,14] = 0.08034633188441791  # S15
$data[13,15] = 0.08034633188441791  # T15
# Row 16
$data[14,0] = 3.0  # E16
$data[14,1] = 1.0  # F16
$data[14,2] = 45.02666966666666  # G16
$data[14,3] = 135.080009  # H16
$data[14,4] = 0.1647240005714903  # I16
$data[14,5] = 0.1647240005714903  # J16
$data[14,6] = 3.0  # K16
$data[14,7] = 1.0  # L16
$data[14,8] = 0.2909853333333334  # M16
$data[14,9] = 0.8729560000000001  # N16
$data[14,10] = 0.03275300561492853  # O16
$data[14,11] = 0.03275300561492853  # P16
$data[14,12] = 13.10210048184489  # Q16
$data[14,13] = 117.918904336604  # R16
$data[14,14] = 0.005395206115631512  # S16
$data[14,15] = 0.005395206115631511  # T16
# Row 17
$data[15,0] = 3.0  # E17
$data[15,1] = 1.0  # F17
$data[15,2] = 45.02666966666666  # G17
$data[15,3] = 135.080009  # H17
$data[15,4] = 0.1647240005714903  # I17
$data[15,5] = 0.1647240005714903  # J17
$data[15,6] = 3.0  # K17
$data[15,7] = 1.0  # L17
$data[15,8] = 0.4814053333333333  # M17
$data[15,9] = 1.444216  # N17
$data[15,10] = 0.0541864822020464  # O17
$data[15,11] = 0.05418648220204641  # P17
$data[15,12] = 21.67607891977155  # Q17
$data[15,13] = 195.084710277944  # R17
$data[15,14] = 0.00892581412521694  # S17
$data[15,15] = 0.008925814125216938  # T17

$ws.Range("E2:T17").Value2 = $data
